$d = $word.ActiveDocument

# Locate the target text in the second ("Phan II") confirmation block:
#   "Hoc ky: 1 ... Nam hoc: 2021-2022"
# and change it to "Hoc ky: 2 ... Nam hoc: 2022-2023", matching the run
# splits produced by a real Word single-character edit (toggling Bold off
# and back on forces the engine to keep the edited character as its own
# run instead of silently re-merging it into a neighbouring run that
# shares the same formatting).

$full = $d.Content.Text

# --- "Hoc ky: 1" -> "Hoc ky: 2" (second occurrence - the filled-in one) ---
$hocKyLabel = "Học kỳ: "
$firstHocKy = $full.IndexOf($hocKyLabel)
$secondHocKy = $full.IndexOf($hocKyLabel, $firstHocKy + 1)
$hocKyValuePos = $secondHocKy + $hocKyLabel.Length

$rHocKy = $d.Range($hocKyValuePos, $hocKyValuePos + 1)
$rHocKy.Font.Bold = $false
$rHocKy.Text = "2"
$rHocKyRestore = $d.Range($hocKyValuePos, $hocKyValuePos + 1)
$rHocKyRestore.Font.Bold = $true

# --- "2021-2022" -> "2022-2023" (second occurrence - the filled-in one) ---
$namHocLabel = "Năm học: "
$firstNamHoc = $full.IndexOf($namHocLabel)
$secondNamHoc = $full.IndexOf($namHocLabel, $firstNamHoc + 1)
$namHocValuePos = $secondNamHoc + $namHocLabel.Length

# "2021-2022"
#  0123456789
# positions that differ from "2022-2023": index 3 (1->2) and index 8 (2->3)
$p1 = $namHocValuePos + 3
$rYear1 = $d.Range($p1, $p1 + 1)
$rYear1.Font.Bold = $false
$rYear1.Text = "2"
$rYear1Restore = $d.Range($p1, $p1 + 1)
$rYear1Restore.Font.Bold = $true

$p2 = $namHocValuePos + 8
$rYear2 = $d.Range($p2, $p2 + 1)
$rYear2.Font.Bold = $false
$rYear2.Text = "3"
$rYear2Restore = $d.Range($p2, $p2 + 1)
$rYear2Restore.Font.Bold = $true

# --- Remove the stray _GoBack bookmark that trailed the year run ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
